# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets of the workbook, per the latest data refresh.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $RowValues) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F$row").Value = $RowValues[$row]
    }
}

# Sheet "展览"
Set-FValues "展览" @{
    2  = 141
    4  = 600
    5  = 2921
    6  = 2921
    7  = 787
    8  = 592
    10 = 410
    11 = 656
    13 = 499
    15 = 2152
    16 = 1259
    17 = 741
    18 = 13
    20 = 2663
    26 = 594
    27 = 594
    28 = 8
    33 = 225
    35 = 388
    36 = 4665
    38 = 14
}

# Sheet "演出"
Set-FValues "演出" @{
    4  = 3
    8  = 330
    9  = 351
    23 = 273
    25 = 296
    27 = 98
    36 = 528
    37 = 528
    38 = 16
}

# Sheet "本地生活"
Set-FValues "本地生活" @{
    4 = 1463
    5 = 568
    6 = 229
    7 = 243
}

# Sheet "全部类型"
Set-FValues "全部类型" @{
    2  = 1463
    3  = 568
    4  = 141
    6  = 229
    9  = 600
    10 = 2921
    11 = 2921
    12 = 787
    13 = 592
    15 = 410
    16 = 656
    17 = 3
    19 = 499
    21 = 330
    22 = 351
    25 = 2152
    26 = 741
    28 = 2663
    33 = 243
    38 = 594
    39 = 594
    40 = 273
    43 = 296
    46 = 388
    48 = 4665
    50 = 528
    51 = 16
}
